$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 corresponds to file "3bcd6935-3054-4dd3-8ac5-6f943fedbf82.md"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 same file -> update Status and Latest Handoff Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-02-25 05:50:18"

# de-de sheet: row 3 same file -> update Status and Latest Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-02-25 05:50:30"
